$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (stock ticker + price) now sourced via Google Finance
$ws.Range("B2").Value = "CHR:TSE"

# Force C2 to stay literal text ("$3.35") instead of being auto-coerced into
# a currency number by the "$" prefix, then drop back to the default style
# so the cell keeps the same (unstyled) formatting it had before the edit.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "`$3.35"
$ws.Range("C2").Style = "Normal"

# Remove the old rows 3-15 that previously held the rest of the scraped data
$ws.Range("A3:C15").EntireRow.Delete()
